# Updated cryptos list with GitHub Actions.
# Note: Price cells that look like plain decimal numbers (e.g. "302.29")
# are written with a leading apostrophe so Excel stores them as text
# (matching the workbook's existing inlineStr cells) instead of silently
# re-typing them as numeric values. Prices that already contain multiple
# "." separators (e.g. "43.056.72") are never auto-converted to numbers,
# so no apostrophe is needed for those.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.056.72"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "2.376.91"
$ws.Range("E3").Value = "  +2.99%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'302.29"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "'97.12"
$ws.Range("E6").Value = "  +1.38%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.501"
$ws.Range("E9").Value = "  +1.82%  "
$ws.Range("D10").Value = "'34.22"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").Value = "'0.0788"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "'0.123"
$ws.Range("E12").Value = "  +2.69%  "
$ws.Range("D13").Value = "'18.37"
$ws.Range("E13").Value = "  -3.57%  "
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").Value = "2.749.85"
$ws.Range("E15").Value = "  +3.06%  "
$ws.Range("D16").Value = "2.380.47"
$ws.Range("E16").Value = "  +3.40%  "
$ws.Range("D17").Value = "'0.806"
$ws.Range("E17").Value = "  +3.00%  "
$ws.Range("D18").Value = "43.040.57"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").Value = "'12.20"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").Value = "'6.30"
$ws.Range("E20").Value = "  +3.98%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "'68.18"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").Value = "'235.48"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'2.25"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'25.00"
$ws.Range("E27").Value = "  +3.00%  "
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").Value = "'9.22"
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("D30").Value = "'31.38"
$ws.Range("E30").Value = "  -2.69%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("D33").Value = "'0.0747"
$ws.Range("E33").Value = "  +6.90%  "
$ws.Range("D34").Value = "'17.41"
$ws.Range("E34").Value = "  -0.89%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.105"
$ws.Range("E35").Value = "  +5.40%  "
$ws.Range("D36").Value = "'1.86"
$ws.Range("E36").Value = "  +6.79%  "
$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").Value = "'23.34"
$ws.Range("E37").Value = "  +17.80%  "
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").Value = "'4.35"
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("D40").Value = "'2.80"
$ws.Range("E40").Value = "  +4.43%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").Value = "'113.95"
$ws.Range("E42").Value = "  -31.06%  "
$ws.Range("D43").Value = "1.948.92"
$ws.Range("D44").Value = "'0.0280"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("E45").Value = "  +2.35%  "
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("D47").Value = "'9.15"
$ws.Range("E47").Value = "  -12.25%  "
$ws.Range("D48").Value = "2.608.16"
$ws.Range("E48").Value = "  +2.81%  "
$ws.Range("E49").Value = "  +3.10%  "
$ws.Range("D50").Value = "'52.45"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'71.86"
$ws.Range("E51").Value = "  -0.17%  "
